$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting for the cells we are about to write, so that values
# like "1.002" or "243.46" are not auto-converted to numeric/date types by Excel,
# matching the inline-string (text) representation in the source workbook.
$cells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12",
    "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17",
    "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "E22",
    "D23", "E23", "E24", "D25", "E25", "B26", "C26", "D26", "E26", "B27",
    "C27", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31",
    "E31", "D32", "E32", "D33", "E33", "D34", "E34", "E35", "D36", "E36",
    "D37", "E37", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42",
    "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47",
    "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($ref in $cells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '29.099.35'
$ws.Range('E2').Value = '  +0.36%  '
$ws.Range('D3').Value = '1.834.95'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.39%  '
$ws.Range('D5').Value = '243.46'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D6').Value = '0.6173'
$ws.Range('E6').Value = '  -2.23%  '
$ws.Range('E7').Value = '  +0.42%  '
$ws.Range('D8').Value = '0.07452'
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('D9').Value = '0.2917'
$ws.Range('E9').Value = '  -0.65%  '
$ws.Range('D10').Value = '23.06'
$ws.Range('E10').Value = '  +0.98%  '
$ws.Range('D11').Value = '0.07698'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '1.838.11'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '4.986'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('D14').Value = '0.6700'
$ws.Range('E14').Value = '  -0.06%  '
$ws.Range('D15').Value = '82.48'
$ws.Range('E15').Value = '  -0.64%  '
$ws.Range('D16').Value = '0.000009285'
$ws.Range('E16').Value = '  -4.37%  '
$ws.Range('D17').Value = '5.925'
$ws.Range('E17').Value = '  -2.58%  '
$ws.Range('D18').Value = '29.086.71'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').Value = '2.094.57'
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').Value = '230.78'
$ws.Range('E20').Value = '  +1.98%  '
$ws.Range('D21').Value = '12.62'
$ws.Range('E21').Value = '  +0.61%  '
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('D23').Value = '7.152'
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = '160.07'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').Value = '8.503'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('B27').Value = 'Stellar'
$ws.Range('C27').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D27').Value = '0.1385'
$ws.Range('E27').Value = '  -1.48%  '
$ws.Range('D28').Value = '17.81'
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('D29').Value = '1.504'
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').Value = '4.155'
$ws.Range('E30').Value = '  +0.88%  '
$ws.Range('D31').Value = '4.117'
$ws.Range('E31').Value = '  +1.14%  '
$ws.Range('D32').Value = '0.05508'
$ws.Range('E32').Value = '  +2.68%  '
$ws.Range('D33').Value = '1.199'
$ws.Range('E33').Value = '  +0.21%  '
$ws.Range('D34').Value = '0.7425'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('E35').Value = '  -1.15%  '
$ws.Range('D36').Value = '1.139'
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('D37').Value = '2.664'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('D39').Value = '1.217.97'
$ws.Range('E39').Value = '  -2.09%  '
$ws.Range('D40').Value = '0.01776'
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('D41').Value = '6.428'
$ws.Range('E41').Value = '  -2.43%  '
$ws.Range('D42').Value = '0.8941'
$ws.Range('D43').Value = '1.002'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').Value = '101.63'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '2.002.55'
$ws.Range('E45').Value = '  +0.85%  '
$ws.Range('D46').Value = '65.48'
$ws.Range('E46').Value = '  +1.20%  '
$ws.Range('D47').Value = '0.00000000122'
$ws.Range('E47').Value = '  -0.75%  '
$ws.Range('D48').Value = '0.5099'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').Value = '0.4061'
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('D50').Value = '9.091'
$ws.Range('E50').Value = '  +1.14%  '
$ws.Range('D51').Value = '0.05833'
$ws.Range('E51').Value = '  +1.20%  '

# Restore default (Normal) style on touched cells so no stray number-format
# style survives the text assignment (keeps cell styling identical to source).
foreach ($ref in $cells) {
    $ws.Range($ref).Style = "Normal"
}

Write-Host "Applied cryptos list update."